$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16511686411348996"
$ws1.Range("B2").Value = "go_stims-16511686411053183.csv"
$ws1.Range("B3").Value = "GNG_stims-16511686411181998.csv"
$ws1.Range("B4").Value = "go_stims-16511686411201637.csv"
$ws1.Range("B5").Value = "GNG_stims-16511686411348996.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16511686461013448"
$ws2.Range("B2").Value = "ZB-match_8-16511686414114819.csv"
$ws2.Range("B3").Value = "TB-16511686433830547.csv"
$ws2.Range("B4").Value = "OB-1651168642279237.csv"
$ws2.Range("B5").Value = "TB-16511686460863428.csv"
$ws2.Range("B6").Value = "OB-1651168641982838.csv"
$ws2.Range("B7").Value = "ZB-match_7-1651168641581166.csv"
$ws2.Range("B8").Value = "OB-16511686417823737.csv"
$ws2.Range("B9").Value = "ZB-match_6-16511686412617369.csv"
$ws2.Range("B10").Value = "TB-16511686446970322.csv"

# --- Sheet 3: RS_TO ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16511686461023467"
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16511686461653464"
$ws4.Range("B2").Value = "MM_stims-16511686461173463.csv"
$ws4.Range("B3").Value = "ZM_stims-1651168646104346.csv"
$ws4.Range("B4").Value = "MM_stims-16511686461483457.csv"
$ws4.Range("B5").Value = "ZM_stims-16511686461183453.csv"
$ws4.Range("B6").Value = "MM_stims-16511686461643467.csv"
$ws4.Range("B7").Value = "ZM_stims-16511686461493454.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16511686462283473"
$ws5.Range("B2").Value = "SAT_stims-16511686461683457.csv"
$ws5.Range("B3").Value = "vSAT_stims-16511686462133448.csv"
$ws5.Range("B4").Value = "vSAT_stims-1651168646196379.csv"
$ws5.Range("B5").Value = "SAT_stims-16511686461803775.csv"
